# Applies the "Update and Spreadsheet input / Preliminary google sheet
# interaction" commit: adds a small Fruit/Lemon sample pair in D9:D10,
# appends a block of new TODO/notes rows (101-110) to Sheet1, adds a new
# external hyperlink for the github compare link, and drops the stale
# "display" text on the existing Google-Sheet hyperlink in A44.
#
# Cell writes are kept in the same order the rows appear in the sheet so
# the shared-string table grows in the same sequence as the source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new rows of notes/TODOs appended after the existing data (row 100) -
$ws.Range("A101").Value = "Crash protection - maybe use the DB?"

# --- small "Fruit / Lemon" sample values next to the existing rows ------
$ws.Range("D9").Value  = "Fruit"
$ws.Range("D10").Value = "Lemon"

$ws.Range("A102").Value = "populate user sheet by role on reset"

# --- A103: new hyperlink row with the github compare link ---------------
$ws.Range("A103").Value = "https://github.com/theoephraim/node-google-spreadsheet/compare/master...reptily:patch-1"
$ws.Hyperlinks.Add($ws.Range("A103"), "https://github.com/theoephraim/node-google-spreadsheet/compare/master...reptily:patch-1")
$ws.Range("A103").Style = "Hyperlink"

$ws.Range("A104").Value = 'FUCK IT! I GIVE UP. Let the admins add more columns and just search for the date keys for inserts'
$ws.Range("A104").Font.Bold = $true

$ws.Range("A105").Value = "Do I need promises"
$ws.Range("A106").Value = "take out dbconnect if not using"
$ws.Range("A107").Value = "5 -> 05 June user input "
$ws.Range("A108").Value = 'for now, im disabling screen population on startup. Debugging takes too long. Also investigate when it sometimes doesn''t populate the dates correctly'
$ws.Range("A109").Value = "The start screen dates seem all our of whack"

# trailing bold marker row (no text, just carries the bold style forward)
$ws.Range("A110").Font.Bold = $true

# --- A44: re-create the existing Google Sheet hyperlink without the
#     stale cached "display" text (URL + " - gid=0") ---------------------
foreach ($hl in @($ws.Hyperlinks)) {
    if ($hl.Range.Row -eq 44) {
        $hl.Delete()
    }
}
$ws.Hyperlinks.Add($ws.Range("A44"), "https://docs.google.com/spreadsheets/d/19zU2Dz78yuttROuU0z54j2S2Fy2gG_gthiE5qrLqsYU/edit", "gid=0")
$ws.Range("A44").Style = "Hyperlink"

# --- restore the selection on the newly-added last row -------------------
$null = $ws.Range("A109").Select()
